$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 572
$ws.Range("I2").Value = 1490
$ws.Range("J2").Value = 6068
$ws.Range("L2").Value = 1621
$ws.Range("M2").Value = 88
$ws.Range("N2").Value = 1056
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 94
$ws.Range("S2").Value = 645
$ws.Range("T2").Value = 1064
$ws.Range("U2").Value = 89
$ws.Range("V2").Value = 9536
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 9438
$ws.Range("Z2").Value = 139
$ws.Range("AA2").Value = 53
